$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 21:40"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6240273
$ws.Range("C4").Value = 24681
$ws.Range("D4").Value = 3469088
$ws.Range("E4").Value = 2582753
$ws.Range("G4").Value = 696
$ws.Range("H4").Value = 188432

# Row 6 - India
$ws.Range("B6").Value = 3766108
$ws.Range("C6").Value = 78169
$ws.Range("D6").Value = 2898087
$ws.Range("E6").Value = 801561
$ws.Range("G6").Value = 1025
$ws.Range("H6").Value = 66460

# Row 12 - España
$ws.Range("B12").Value = 470973
$ws.Range("C12").Value = 8115

# Row 23 - Alemania
$ws.Range("B23").Value = 245957
$ws.Range("C23").Value = 1165
$ws.Range("E23").Value = 16677
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 9380

# Row 61 - Suiza
$ws.Range("D61").Value = 36300
$ws.Range("E61").Value = 4082

# Row 95 - Guinea
$ws.Range("B95").Value = 9479
$ws.Range("C95").Value = 70
$ws.Range("D95").Value = 8527
$ws.Range("E95").Value = 893

# Row 104 - Mauritania
$ws.Range("B104").Value = 7075
$ws.Range("C104").Value = 27
$ws.Range("E104").Value = 452

# Row 110 - Guinea Ecuatorial
$ws.Range("B110").Value = 4965
$ws.Range("C110").Value = 24
$ws.Range("E110").Value = 998
